$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4799.8
$ws.Range("I69").Value = 2999.5
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 8998.5
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -8124.5
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 4799.8
$ws.Range("I72").Value = 2999.5
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 26995.5
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -22627.5
$ws.Range("N72").Value = -62736
$ws.Range("H74").Value = 4325.1177
$ws.Range("I74").Value = 2440.875
$ws.Range("K74").Value = 2440.875
$ws.Range("M74").Value = -1504.875
$ws.Range("H77").Value = 4325.1177
$ws.Range("I77").Value = 2440.875
$ws.Range("K77").Value = 12204.375
$ws.Range("M77").Value = -7524.375
$ws.Range("H86").Value = 2377.5312
$ws.Range("I86").Value = 1603.2142
$ws.Range("J86").Value = 2979.7778
$ws.Range("K86").Value = 1603.2142
$ws.Range("L86").Value = 2979.7778
$ws.Range("M86").Value = -480.2141999999999
$ws.Range("N86").Value = -5225.7778
$ws.Range("H87").Value = 61070
$ws.Range("J87").Value = 61070
$ws.Range("L87").Value = 61070
$ws.Range("N87").Value = -63566
$ws.Range("H89").Value = 2377.5312
$ws.Range("I89").Value = 1603.2142
$ws.Range("J89").Value = 2979.7778
$ws.Range("K89").Value = 8016.071
$ws.Range("L89").Value = 14898.889
$ws.Range("M89").Value = -2400.071
$ws.Range("N89").Value = -26130.889
$ws.Range("H90").Value = 61070
$ws.Range("J90").Value = 61070
$ws.Range("L90").Value = 183210
$ws.Range("N90").Value = -195690
$ws.Range("H96").Value = 2587.2307
$ws.Range("I96").Value = 1514.8889
$ws.Range("K96").Value = 4544.6667
$ws.Range("M96").Value = -3171.6667
$ws.Range("H100").Value = 2988
$ws.Range("I100").Value = 2334.1667
$ws.Range("J100").Value = 4949.5
$ws.Range("K100").Value = 2334.1667
$ws.Range("L100").Value = 4949.5
$ws.Range("M100").Value = -1793.1667
$ws.Range("N100").Value = -6031.5
$ws.Range("H107").Value = 539.7083
$ws.Range("I107").Value = 557.9474
$ws.Range("J107").Value = 470.4
$ws.Range("K107").Value = 557.9474
$ws.Range("L107").Value = 470.4
$ws.Range("M107").Value = 1362.0526
$ws.Range("N107").Value = -4310.4
$ws.Range("H118").Value = 1130.2
$ws.Range("I118").Value = 1063.5834
$ws.Range("J118").Value = 1396.6666
$ws.Range("K118").Value = 3190.7502
$ws.Range("L118").Value = 4189.9998
$ws.Range("M118").Value = -1533.7502
$ws.Range("N118").Value = -7503.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14398.464
$ws.Range("I32").Value = 11959.818
$ws.Range("K32").Value = 11959.818
$ws.Range("M32").Value = -11672.818
$ws.Range("H61").Value = 14290434
$ws.Range("I61").Value = 20003292
$ws.Range("K61").Value = 20003292
$ws.Range("M61").Value = -20003080
$ws.Range("H74").Value = 26346818
$ws.Range("I74").Value = 32295506
$ws.Range("J74").Value = 2629.4285
$ws.Range("K74").Value = 32295506
$ws.Range("L74").Value = 2629.4285
$ws.Range("M74").Value = -32294632
$ws.Range("N74").Value = -4377.4285
$ws.Range("H77").Value = 26346818
$ws.Range("I77").Value = 32295506
$ws.Range("J77").Value = 2629.4285
$ws.Range("K77").Value = 161477530
$ws.Range("L77").Value = 13147.1425
$ws.Range("M77").Value = -161473162
$ws.Range("N77").Value = -21883.1425
$ws.Range("H110").Value = 20509.21
$ws.Range("I110").Value = 24779
$ws.Range("K110").Value = 24779
$ws.Range("M110").Value = -22734
$ws.Range("H136").Value = 14290434
$ws.Range("I136").Value = 20003292
$ws.Range("K136").Value = 60009876
$ws.Range("M136").Value = -60007326

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 19626.107
$ws.Range("I86").Value = 10048.292
$ws.Range("K86").Value = 10048.292
$ws.Range("M86").Value = -8925.291999999999
$ws.Range("H89").Value = 19626.107
$ws.Range("I89").Value = 10048.292
$ws.Range("K89").Value = 50241.46
$ws.Range("M89").Value = -44625.46
$ws.Range("H94").Value = 1438.0667
$ws.Range("I94").Value = 283.75
$ws.Range("J94").Value = 1857.8182
$ws.Range("K94").Value = 283.75
$ws.Range("L94").Value = 1857.8182
$ws.Range("M94").Value = 167.25
$ws.Range("N94").Value = -2759.8182
$ws.Range("H99").Value = 4024.3076
$ws.Range("I99").Value = 3230.7222
$ws.Range("K99").Value = 3230.7222
$ws.Range("M99").Value = -1732.7222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1200
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -913
$ws.Range("H22").Value = 8718.416999999999
$ws.Range("I22").Value = 8718.416999999999
$ws.Range("K22").Value = 8718.416999999999
$ws.Range("M22").Value = -8368.416999999999
$ws.Range("H31").Value = 4054.818
$ws.Range("I31").Value = 3214.9333
$ws.Range("K31").Value = 3214.9333
$ws.Range("M31").Value = -2919.9333
$ws.Range("H34").Value = 4054.818
$ws.Range("I34").Value = 3214.9333
$ws.Range("K34").Value = 3214.9333
$ws.Range("M34").Value = -3012.9333
$ws.Range("H105").Value = 10407.75
$ws.Range("I105").Value = 1231.7142
$ws.Range("J105").Value = 23254.2
$ws.Range("K105").Value = 1231.7142
$ws.Range("L105").Value = 23254.2
$ws.Range("M105").Value = 515.2858000000001
$ws.Range("N105").Value = -26748.2
$ws.Range("H113").Value = 1200
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970
$ws.Range("H132").Value = 56898.055
$ws.Range("J132").Value = 4310.9165
$ws.Range("L132").Value = 12932.7495
$ws.Range("N132").Value = -17992.7495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 159303.36
$ws.Range("J37").Value = 159303.36
$ws.Range("L37").Value = 477910.08
$ws.Range("N37").Value = -478134.08
$ws.Range("H62").Value = 6670928
$ws.Range("I62").Value = 4950
$ws.Range("K62").Value = 14850
$ws.Range("M62").Value = -14164
$ws.Range("H64").Value = 15206.25
$ws.Range("I64").Value = 3075
$ws.Range("J64").Value = 19250
$ws.Range("K64").Value = 9225
$ws.Range("L64").Value = 57750
$ws.Range("M64").Value = -8955
$ws.Range("N64").Value = -58290
$ws.Range("H65").Value = 6670928
$ws.Range("I65").Value = 4950
$ws.Range("K65").Value = 44550
$ws.Range("M65").Value = -41118
$ws.Range("H67").Value = 15206.25
$ws.Range("I67").Value = 3075
$ws.Range("J67").Value = 19250
$ws.Range("K67").Value = 9225
$ws.Range("L67").Value = 57750
$ws.Range("M67").Value = -8289
$ws.Range("N67").Value = -59622
$ws.Range("H70").Value = 4456.2856
$ws.Range("I70").Value = 324.25
$ws.Range("J70").Value = 9965.666999999999
$ws.Range("K70").Value = 972.75
$ws.Range("L70").Value = 29897.001
$ws.Range("M70").Value = -657.75
$ws.Range("N70").Value = -30527.001
$ws.Range("H73").Value = 4456.2856
$ws.Range("I73").Value = 324.25
$ws.Range("J73").Value = 9965.666999999999
$ws.Range("K73").Value = 972.75
$ws.Range("L73").Value = 29897.001
$ws.Range("M73").Value = 119.25
$ws.Range("N73").Value = -32081.001
$ws.Range("H122").Value = 1700.3
$ws.Range("J122").Value = 2143.4285
$ws.Range("L122").Value = 19290.8565
$ws.Range("N122").Value = -24190.8565
$ws.Range("H129").Value = 1722.8387
$ws.Range("I129").Value = 551
$ws.Range("J129").Value = 2130.4348
$ws.Range("K129").Value = 1653
$ws.Range("L129").Value = 6391.3044
$ws.Range("M129").Value = 3347
$ws.Range("N129").Value = -16391.3044
$ws.Range("H140").Value = 1014
$ws.Range("I140").Value = 917.9048
$ws.Range("K140").Value = 2753.7144
$ws.Range("M140").Value = 2426.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2902.1765
$ws.Range("J122").Value = 2912.25
$ws.Range("L122").Value = 8736.75
$ws.Range("N122").Value = -13636.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 679.0454999999999
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 500
$ws.Range("M55").Value = -327
$ws.Range("H93").Value = 2211.1177
$ws.Range("I93").Value = 1123.4
$ws.Range("J93").Value = 2664.3333
$ws.Range("K93").Value = 1123.4
$ws.Range("L93").Value = 2664.3333
$ws.Range("M93").Value = 124.5999999999999
$ws.Range("N93").Value = -5160.3333
$ws.Range("H136").Value = 1886.8158
$ws.Range("I136").Value = 668.8889
$ws.Range("K136").Value = 2006.6667
$ws.Range("M136").Value = 543.3332999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1177.3846
$ws.Range("I81").Value = 1186.5
$ws.Range("J81").Value = 1162.8
$ws.Range("K81").Value = 2373
$ws.Range("L81").Value = 2325.6
$ws.Range("M81").Value = -1312
$ws.Range("N81").Value = -4447.6
$ws.Range("H84").Value = 1177.3846
$ws.Range("I84").Value = 1186.5
$ws.Range("J84").Value = 1162.8
$ws.Range("K84").Value = 11865
$ws.Range("L84").Value = 11628
$ws.Range("M84").Value = -6561
$ws.Range("N84").Value = -22236
$ws.Range("H103").Value = 15313.556
$ws.Range("J103").Value = 15313.556
$ws.Range("L103").Value = 15313.556
$ws.Range("N103").Value = -17657.556
$ws.Range("H106").Value = 36678.445
$ws.Range("I106").Value = 24000
$ws.Range("J106").Value = 43017.668
$ws.Range("K106").Value = 24000
$ws.Range("L106").Value = 43017.668
$ws.Range("M106").Value = -22738
